# Apply text styling updates across all slides:
#  - Title placeholders: explicitly set Italic = False (keeps existing Bold)
#  - Body/content placeholders: explicitly set Bold = False and Italic = False
$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)

    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $shp = $s.Shapes.Item($j)

        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange

            if ($shp.PlaceholderFormat.Type -eq 1) {
                # Title placeholder
                $tr.Font.Italic = 0
            } else {
                # Body / content placeholder
                $tr.Font.Bold = 0
                $tr.Font.Italic = 0
            }
        }
    }
}
